$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,13
$data[0,0] = 0.7667715243641169
$data[0,1] = 0.2713778478185418
$data[0,2] = 0.08107532069041667
$data[0,3] = 0.1350018207780934
$data[0,4] = 2.729558823629873
$data[0,5] = 0
$data[0,6] = 0.07973214163530429
$data[0,7] = 1.369635133289769
$data[0,8] = 0.2156021860240713
$data[0,9] = 1.030382707340493
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 2.548320960721085
$data[1,0] = 0.7263890224817544
$data[1,1] = 0.2603134060759658
$data[1,2] = 0.07855395878265625
$data[1,3] = 0.1318959325522506
$data[1,4] = 2.71753019567916
$data[1,5] = 0
$data[1,6] = 0.07973214163530429
$data[1,7] = 1.369253655889828
$data[1,8] = 0.2113600284143686
$data[1,9] = 0.9797632597493759
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 2.564830886544563
$data[2,0] = 0.7020253158227376
$data[2,1] = 0.2536840887005667
$data[2,2] = 0.07704666928802339
$data[2,3] = 0.1300628855400099
$data[2,4] = 2.711447463746978
$data[2,5] = 0
$data[2,6] = 0.07973214163530429
$data[2,7] = 1.369619534222117
$data[2,8] = 0.2088774292417952
$data[2,9] = 0.9492766824184287
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 2.575691282604808
$data[3,0] = 0.6922053995450312
$data[3,1] = 0.2510237924602166
$data[3,2] = 0.07644274143563479
$data[3,3] = 0.1293345034873283
$data[3,4] = 2.709296048314144
$data[3,5] = 0
$data[3,6] = 0.07973214163530429
$data[3,7] = 1.369919395441826
$data[3,8] = 0.2078964048938019
$data[3,9] = 0.9370024204895344
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 2.58029867709061
$data[4,0] = 0.690581363775209
$data[4,1] = 0.2505845378130402
$data[4,2] = 0.07634308313520677
$data[4,3] = 0.1292146794776379
$data[4,4] = 2.70895857048167
$data[4,5] = 0
$data[4,6] = 0.07973214163530429
$data[4,7] = 1.369978288872019
$data[4,8] = 0.2077353566308275
$data[4,9] = 0.9349733014612411
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 2.581074700493801
$data[5,0] = 0.7018924416694858
$data[5,1] = 0.2536480443219489
$data[5,2] = 0.07703848273320801
$data[5,3] = 0.1300529870200577
$data[5,4] = 2.711417123964551
$data[5,5] = 0
$data[5,6] = 0.07973214163530429
$data[5,7] = 1.369622968012912
$data[5,8] = 0.2088640747467352
$data[5,9] = 0.9491105431762605
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 2.575752684098731
$data[6,0] = 0.7527581237681318
$data[6,1] = 0.2675286589347934
$data[6,2] = 0.08019749994583236
$data[6,3] = 0.1339155521334803
$data[6,4] = 2.725140758814689
$data[6,5] = 0
$data[6,6] = 0.07973214163530429
$data[6,7] = 1.369378960174032
$data[6,8] = 0.214114116242186
$data[6,9] = 1.012805768017472
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 2.553863374647477
$data[7,0] = 0.8559326802435123
$data[7,1] = 0.296058816568916
$data[7,2] = 0.08671518670705325
$data[7,3] = 0.1420777950149485
$data[7,4] = 2.762407981765506
$data[7,5] = 0
$data[7,6] = 0.07973214163530429
$data[7,7] = 1.373669887941467
$data[7,8] = 0.2253812985054395
$data[7,9] = 1.142437021404817
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 2.516684461755077
$data[8,0] = 0.9338402286480232
$data[8,1] = 0.3178306191757656
$data[8,2] = 0.09169963031988004
$data[8,3] = 0.1484348237307174
$data[8,4] = 2.796132083912909
$data[8,5] = 0
$data[8,6] = 0.07973214163530429
$data[8,7] = 1.379743150831587
$data[8,8] = 0.2342571856878664
$data[8,9] = 1.240586262868874
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 2.492882011755839
$data[9,0] = 0.9697435849798808
$data[9,1] = 0.3279139580942001
$data[9,2] = 0.09400957029909307
$data[9,3] = 0.1514054942527068
$data[9,4] = 2.812858823555274
$data[9,5] = 0
$data[9,6] = 0.07973214163530429
$data[9,7] = 1.383143276038645
$data[9,8] = 0.2384261352473231
$data[9,9] = 1.28587555827599
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 2.48281869244326
$data[10,0] = 0.9834059455105262
$data[10,1] = 0.3317582010884905
$data[10,2] = 0.09489036937669937
$data[10,3] = 0.1525417665049886
$data[10,4] = 2.8193925047389
$data[10,5] = 0
$data[10,6] = 0.07973214163530429
$data[10,7] = 1.384522680477382
$data[10,8] = 0.2400237617332834
$data[10,9] = 1.3031178968738
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 2.479118127891965
$data[11,0] = 0.9804605525687293
$data[11,1] = 0.3309291205552256
$data[11,2] = 0.09470040390853285
$data[11,3] = 0.1522965452825318
$data[11,4] = 2.817976474662643
$data[11,5] = 0
$data[11,6] = 0.07973214163530429
$data[11,7] = 1.384221512768022
$data[11,8] = 0.2396788411433022
$data[11,9] = 1.29940034505529
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 2.479910205533386
$data[12,0] = 0.9708662621333133
$data[12,1] = 0.3282297064401121
$data[12,2] = 0.09408191261927357
$data[12,3] = 0.15149874865633
$data[12,4] = 2.813392350290613
$data[12,5] = 0
$data[12,6] = 0.07973214163530429
$data[12,7] = 1.38325491849843
$data[12,8] = 0.2385571930909407
$data[12,9] = 1.287292245936527
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 2.482512034789238
$data[13,0] = 0.9649981471704905
$data[13,1] = 0.3265796142914894
$data[13,2] = 0.09370385894267486
$data[13,3] = 0.1510115527514202
$data[13,4] = 2.810610451141017
$data[13,5] = 0
$data[13,6] = 0.07973214163530429
$data[13,7] = 1.382674819222736
$data[13,8] = 0.2378726190697336
$data[13,9] = 1.279887705305896
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 2.484120088820973
$data[14,0] = 0.9315031716757858
$data[14,1] = 0.3171752668032184
$data[14,2] = 0.09154952232107405
$data[14,3] = 0.148242270340873
$data[14,4] = 2.795066865316443
$data[14,5] = 0
$data[14,6] = 0.07973214163530429
$data[14,7] = 1.379533786977774
$data[14,8] = 0.2339873815184319
$data[14,9] = 1.237639404777156
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 2.493555079582023
$data[15,0] = 0.9110736169458278
$data[15,1] = 0.3114520009980879
$data[15,2] = 0.09023876444084067
$data[15,3] = 0.1465636005731312
$data[15,4] = 2.785886513406737
$data[15,5] = 0
$data[15,6] = 0.07973214163530429
$data[15,7] = 1.377770250478875
$data[15,8] = 0.2316375764879268
$data[15,9] = 1.211885634241668
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 2.499539137773297
$data[16,0] = 0.899366640933863
$data[16,1] = 0.3081769927461266
$data[16,2] = 0.08948885277109753
$data[16,3] = 0.1456054930735107
$data[16,4] = 2.78073660107033
$data[16,5] = 0
$data[16,6] = 0.07973214163530429
$data[16,7] = 1.376815890890114
$data[16,8] = 0.2302983823198161
$data[16,9] = 1.197133045249188
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 2.503052955670555
$data[17,0] = 0.8954103426145252
$data[17,1] = 0.3070710235913907
$data[17,2] = 0.08923563392454525
$data[17,3] = 0.1452823682927118
$data[17,4] = 2.779015310707862
$data[17,5] = 0
$data[17,6] = 0.07973214163530429
$data[17,7] = 1.376503056773842
$data[17,8] = 0.2298470737540725
$data[17,9] = 1.192148424559775
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 2.504255023080361
$data[18,0] = 0.9132438693823417
$data[18,1] = 0.3120595062557356
$data[18,2] = 0.09037788296335236
$data[18,3] = 0.1467415299417496
$data[18,4] = 2.786850280489276
$data[18,5] = 0
$data[18,6] = 0.07973214163530429
$data[18,7] = 1.37795177281636
$data[18,8] = 0.2318864384756978
$data[18,9] = 1.214620925465169
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 2.498894676571567
$data[19,0] = 0.9736825320529761
$data[19,1] = 0.3290218858192873
$data[19,2] = 0.09426341391038306
$data[19,3] = 0.1517327729250368
$data[19,4] = 2.814733397683355
$data[19,5] = 0
$data[19,6] = 0.07973214163530429
$data[19,7] = 1.383536336264996
$data[19,8] = 0.2388861339178874
$data[19,9] = 1.290846183014111
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 2.481744821705547
$data[20,0] = 1.013570527292586
$data[20,1] = 0.3402588151399186
$data[20,2] = 0.09683822632078432
$data[20,3] = 0.1550609635158651
$data[20,4] = 2.834120401188855
$data[20,5] = 0
$data[20,6] = 0.07973214163530429
$data[20,7] = 1.387721642510087
$data[20,8] = 0.2435712402599393
$data[20,9] = 1.34120165432185
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 2.471178829844646
$data[21,0] = 0.9922460700461784
$data[21,1] = 0.3342475922457595
$data[21,2] = 0.09546077353677163
$data[21,3] = 0.1532785921243374
$data[21,4] = 2.823666573541075
$data[21,5] = 0
$data[21,6] = 0.07973214163530429
$data[21,7] = 1.385438803817422
$data[21,8] = 0.2410605890286348
$data[21,9] = 1.314276732149125
$data[21,10] = 0
$data[21,11] = 0
$data[21,12] = 2.476759230073668
$data[22,0] = 0.9122625791182202
$data[22,1] = 0.3117848052527563
$data[22,2] = 0.09031497607631422
$data[22,3] = 0.1466610663279013
$data[22,4] = 2.786414162706535
$data[22,5] = 0
$data[22,6] = 0.07973214163530429
$data[22,7] = 1.377869521164669
$data[22,8] = 0.2317738914090484
$data[22,9] = 1.213384133365366
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 2.499185808514724
$data[23,0] = 0.8276527040742963
$data[23,1] = 0.2881991530480832
$data[23,2] = 0.08491751243845158
$data[23,3] = 0.1398065722140842
$data[23,4] = 2.751214431737068
$data[23,5] = 0
$data[23,6] = 0.07973214163530429
$data[23,7] = 1.371997131505609
$data[23,8] = 0.2222286259927131
$data[23,9] = 1.106859545683704
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 2.526126182015986

$ws.Range("B2:N25").Value = $data
Write-Host "Applied updated values to B2:N25"
